$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1670.875
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 1689.2903
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 5067.8709
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -7283.8709

$ws.Range("H129").Value = 1040.7778
$ws.Range("I129").Value = 660
$ws.Range("J129").Value = 1088.375
$ws.Range("K129").Value = 1980
$ws.Range("L129").Value = 3265.125
$ws.Range("M129").Value = 3020
$ws.Range("N129").Value = -13265.125

$ws.Range("H138").Value = 2273.94
$ws.Range("I138").Value = 2195.4614
$ws.Range("J138").Value = 2285.6667
$ws.Range("K138").Value = 6586.3842
$ws.Range("L138").Value = 6857.000100000001
$ws.Range("M138").Value = -1446.3842
$ws.Range("N138").Value = -17137.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12910.272
$ws.Range("I32").Value = 12820.951
$ws.Range("J32").Value = 14000
$ws.Range("K32").Value = 12820.951
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = -12533.951
$ws.Range("N32").Value = -14574

$ws.Range("H61").Value = 1370.0588
$ws.Range("I61").Value = 1330.6875
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1330.6875
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1118.6875
$ws.Range("N61").Value = -2424

$ws.Range("H63").Value = 5503.222
$ws.Range("I63").Value = 4941.125
$ws.Range("K63").Value = 4941.125
$ws.Range("M63").Value = -4255.125

$ws.Range("H66").Value = 5503.222
$ws.Range("I66").Value = 4941.125
$ws.Range("K66").Value = 24705.625
$ws.Range("M66").Value = -21273.625

$ws.Range("H132").Value = 1521.4412
$ws.Range("I132").Value = 1038.72
$ws.Range("J132").Value = 2862.3333
$ws.Range("K132").Value = 3116.16
$ws.Range("L132").Value = 8586.999899999999
$ws.Range("M132").Value = -586.1599999999999
$ws.Range("N132").Value = -13646.9999

$ws.Range("H136").Value = 1370.0588
$ws.Range("I136").Value = 1330.6875
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 3992.0625
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -1442.0625
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 70084.125
$ws.Range("J132").Value = 70084.125
$ws.Range("L132").Value = 70084.125
$ws.Range("N132").Value = -80204.125

$ws.Range("H134").Value = 24392236
$ws.Range("I134").Value = 26317518
$ws.Range("J134").Value = 5331.3335
$ws.Range("K134").Value = 78952554
$ws.Range("L134").Value = 15994.0005
$ws.Range("M134").Value = -78950019
$ws.Range("N134").Value = -21064.0005

$ws.Range("H135").Value = 62760
$ws.Range("J135").Value = 62760
$ws.Range("L135").Value = 62760
$ws.Range("N135").Value = -72900

$ws.Range("H137").Value = 63570
$ws.Range("J137").Value = 63570
$ws.Range("L137").Value = 63570
$ws.Range("N137").Value = -73770

$ws.Range("H138").Value = 57790
$ws.Range("J138").Value = 57790
$ws.Range("L138").Value = 57790
$ws.Range("N138").Value = -68070

$ws.Range("H140").Value = 52574
$ws.Range("J140").Value = 52574
$ws.Range("L140").Value = 52574
$ws.Range("N140").Value = -62934

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 74274.42999999999
$ws.Range("I62").Value = 102624.2
$ws.Range("J62").Value = 3400
$ws.Range("K62").Value = 102624.2
$ws.Range("L62").Value = 3400
$ws.Range("M62").Value = -102000.2
$ws.Range("N62").Value = -4648

$ws.Range("H65").Value = 74274.42999999999
$ws.Range("I65").Value = 102624.2
$ws.Range("J65").Value = 3400
$ws.Range("K65").Value = 513121
$ws.Range("L65").Value = 17000
$ws.Range("M65").Value = -510001
$ws.Range("N65").Value = -23240

$ws.Range("H134").Value = 1295.2106
$ws.Range("I134").Value = 1322.7778
$ws.Range("J134").Value = 799
$ws.Range("K134").Value = 3968.3334
$ws.Range("L134").Value = 2397
$ws.Range("M134").Value = -1433.3334
$ws.Range("N134").Value = -7467

$ws.Range("H138").Value = 41290
$ws.Range("J138").Value = 41290
$ws.Range("L138").Value = 41290
$ws.Range("N138").Value = -51570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2724.8667
$ws.Range("I122").Value = 2782.68
$ws.Range("J122").Value = 2435.8
$ws.Range("K122").Value = 8348.039999999999
$ws.Range("L122").Value = 7307.400000000001
$ws.Range("M122").Value = -5898.039999999999
$ws.Range("N122").Value = -12207.4

$ws.Range("H126").Value = 2631.8667
$ws.Range("I126").Value = 2093.6
$ws.Range("J126").Value = 3708.4
$ws.Range("K126").Value = 6280.799999999999
$ws.Range("L126").Value = 11125.2
$ws.Range("M126").Value = -3810.799999999999
$ws.Range("N126").Value = -16065.2

$ws.Range("H132").Value = 2114.7576
$ws.Range("I132").Value = 1246.35
$ws.Range("J132").Value = 3450.7693
$ws.Range("K132").Value = 3739.05
$ws.Range("L132").Value = 10352.3079
$ws.Range("M132").Value = -1209.05
$ws.Range("N132").Value = -15412.3079

$ws.Range("H135").Value = 47172.855
$ws.Range("J135").Value = 47172.855
$ws.Range("L135").Value = 47172.855
$ws.Range("N135").Value = -57312.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 70715
$ws.Range("J127").Value = 70715
$ws.Range("L127").Value = 70715
$ws.Range("N127").Value = -80635

$ws.Range("H132").Value = 5267.3945
$ws.Range("I132").Value = 5133.1875
$ws.Range("K132").Value = 15399.5625
$ws.Range("M132").Value = -12869.5625

$ws.Range("H133").Value = 54181.2
$ws.Range("J133").Value = 56835.11
$ws.Range("L133").Value = 56835.11
$ws.Range("N133").Value = -61895.11

$ws.Range("H136").Value = 1669.2565
$ws.Range("I136").Value = 1684.2368
$ws.Range("J136").Value = 1100
$ws.Range("K136").Value = 5052.7104
$ws.Range("L136").Value = 3300
$ws.Range("M136").Value = -2502.7104
$ws.Range("N136").Value = -8400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4930.3
$ws.Range("I122").Value = 3180.6
$ws.Range("J122").Value = 6680
$ws.Range("K122").Value = 9541.799999999999
$ws.Range("L122").Value = 20040
$ws.Range("M122").Value = -7091.799999999999
$ws.Range("N122").Value = -24940

$ws.Range("H132").Value = 1214.4043
$ws.Range("I132").Value = 1058.591
$ws.Range("J132").Value = 3499.6667
$ws.Range("K132").Value = 3175.773
$ws.Range("L132").Value = 10499.0001
$ws.Range("M132").Value = -645.7729999999997
$ws.Range("N132").Value = -15559.0001

$ws.Range("H136").Value = 816.1842
$ws.Range("I136").Value = 806.3333
$ws.Range("J136").Value = 840.36365
$ws.Range("K136").Value = 2418.9999
$ws.Range("L136").Value = 2521.09095
$ws.Range("M136").Value = 131.0001000000002
$ws.Range("N136").Value = -7621.09095
